$d = $word.ActiveDocument

# 1) Programa: insert two manual line breaks inside the long paragraph
$d.Content.Find.Execute("06 Função das proteínasMioglobina", $true, $false, $false, $false, $false,
                         $true, 1, $false, "06 Função das proteínas^lMioglobina", 2)

$d.Content.Find.Execute("catálise enzimáticaNomenclatura das enzimas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "catálise enzimática^lNomenclatura das enzimas", 2)

# 2) Avaliação / Critério: insert a manual line break between the trailing
#    spaces after "disciplina." and "A Nota final"
$d.Content.Find.Execute("disciplina.            A Nota final", $true, $false, $false, $false, $false,
                         $true, 1, $false, "disciplina.            ^lA Nota final", 2)

# 3) Bibliografia: insert a manual line break before "2. Voet"
$d.Content.Find.Execute("7ª Ed. 20192. Voet", $true, $false, $false, $false, $false,
                         $true, 1, $false, "7ª Ed. 2019^l2. Voet", 2)
